{"js": "// Update each \"two-digit \u00f7 one-digit =\" expression in the worksheet table\n// to the new value from the target revision. Each original expression is\n// unique in the document, so we can safely search for it and replace the\n// whole match with the new text (replacing the run's text, so any run\n// formatting such as font/size is preserved).\nconst pairs = [\n  [\"50\u00f73=\", \"89\u00f73=\"],\n  [\"53\u00f72=\", \"43\u00f79=\"],\n  [\"26\u00f74=\", \"35\u00f78=\"],\n  [\"56\u00f74=\", \"49\u00f76=\"],\n  [\"11\u00f79=\", \"86\u00f79=\"],\n  [\"84\u00f72=\", \"62\u00f78=\"],\n  [\"65\u00f73=\", \"21\u00f73=\"],\n  [\"92\u00f79=\", \"66\u00f73=\"],\n  [\"28\u00f72=\", \"37\u00f75=\"],\n  [\"91\u00f73=\", \"47\u00f75=\"],\n  [\"86\u00f73=\", \"57\u00f79=\"],\n  [\"53\u00f79=\", \"25\u00f79=\"],\n  [\"55\u00f74=\", \"95\u00f77=\"],\n  [\"44\u00f79=\", \"51\u00f73=\"],\n  [\"39\u00f79=\", \"29\u00f79=\"],\n  [\"87\u00f79=\", \"56\u00f77=\"],\n  [\"43\u00f73=\", \"28\u00f75=\"],\n  [\"98\u00f76=\", \"33\u00f74=\"],\n  [\"29\u00f74=\", \"80\u00f77=\"],\n  [\"88\u00f77=\", \"88\u00f73=\"],\n  [\"50\u00f77=\", \"26\u00f78=\"],\n  [\"44\u00f73=\", \"84\u00f75=\"],\n  [\"60\u00f78=\", \"23\u00f77=\"],\n  [\"71\u00f79=\", \"60\u00f77=\"],\n  [\"40\u00f78=\", \"77\u00f74=\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update each \"two-digit \u00f7 one-digit =\" expression in the worksheet table to\n# the new value from the target revision. Each original expression occurs\n# exactly once in the document, so a plain Find/Replace (wdReplaceAll is\n# used defensively, but each pattern only ever matches a single run) swaps\n# the old text for the new text while leaving all run/paragraph formatting\n# (font, size, table layout, etc.) untouched.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n\n$find.Execute(\"50\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"89\u00f73=\", 2)\n$find.Execute(\"53\u00f72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"43\u00f79=\", 2)\n$find.Execute(\"26\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"35\u00f78=\", 2)\n$find.Execute(\"56\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00f76=\", 2)\n$find.Execute(\"11\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"86\u00f79=\", 2)\n$find.Execute(\"84\u00f72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"62\u00f78=\", 2)\n$find.Execute(\"65\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"21\u00f73=\", 2)\n$find.Execute(\"92\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"66\u00f73=\", 2)\n$find.Execute(\"28\u00f72=\", $false, $false, $false, $false, $false, $true, 1, $false, \"37\u00f75=\", 2)\n$find.Execute(\"91\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"47\u00f75=\", 2)\n$find.Execute(\"86\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"57\u00f79=\", 2)\n$find.Execute(\"53\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"25\u00f79=\", 2)\n$find.Execute(\"55\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"95\u00f77=\", 2)\n$find.Execute(\"44\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"51\u00f73=\", 2)\n$find.Execute(\"39\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"29\u00f79=\", 2)\n$find.Execute(\"87\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00f77=\", 2)\n$find.Execute(\"43\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"28\u00f75=\", 2)\n$find.Execute(\"98\u00f76=\", $false, $false, $false, $false, $false, $true, 1, $false, \"33\u00f74=\", 2)\n$find.Execute(\"29\u00f74=\", $false, $false, $false, $false, $false, $true, 1, $false, \"80\u00f77=\", 2)\n$find.Execute(\"88\u00f77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00f73=\", 2)\n$find.Execute(\"50\u00f77=\", $false, $false, $false, $false, $false, $true, 1, $false, \"26\u00f78=\", 2)\n$find.Execute(\"44\u00f73=\", $false, $false, $false, $false, $false, $true, 1, $false, \"84\u00f75=\", 2)\n$find.Execute(\"60\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"23\u00f77=\", 2)\n$find.Execute(\"71\u00f79=\", $false, $false, $false, $false, $false, $true, 1, $false, \"60\u00f77=\", 2)\n$find.Execute(\"40\u00f78=\", $false, $false, $false, $false, $false, $true, 1, $false, \"77\u00f74=\", 2)\n"}
